$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal string value to a cell without Excel's
# automatic number/date coercion overwriting the text (e.g. "165.10" -> 165.1).
# Temporarily flips the cell to Text format, writes the literal value, then
# restores the cell's original style so no visible formatting changes.
function Set-LiteralText {
    param($Range, [string]$Text)
    $origStyle = $Range.Style
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = $origStyle
}

Set-LiteralText $ws.Range("D2") '63.149.51'
$ws.Range("E2").Value = '  -2.11%  '
Set-LiteralText $ws.Range("D3") '2.570.19'
$ws.Range("E3").Value = '  -3.51%  '
$ws.Range("E4").Value = '  +0.08%  '
Set-LiteralText $ws.Range("D5") '586.19'
$ws.Range("E5").Value = '  -3.68%  '
Set-LiteralText $ws.Range("D6") '148.54'
$ws.Range("E6").Value = '  -2.66%  '
$ws.Range("E7").Value = '  +0.05%  '
Set-LiteralText $ws.Range("D8") '0.584'
$ws.Range("E8").Value = '  -1.27%  '
$ws.Range("E9").Value = '  -1.14%  '
Set-LiteralText $ws.Range("D10") '5.74'
$ws.Range("E10").Value = '  +2.20%  '
Set-LiteralText $ws.Range("D11") '0.381'
$ws.Range("E11").Value = '  -1.90%  '
$ws.Range("E12").Value = '  -0.88%  '
Set-LiteralText $ws.Range("D13") '27.35'
$ws.Range("E13").Value = '  -3.07%  '
Set-LiteralText $ws.Range("D14") '3.033.24'
$ws.Range("E14").Value = '  -3.48%  '
Set-LiteralText $ws.Range("D15") '63.064.80'
$ws.Range("E15").Value = '  -2.00%  '
Set-LiteralText $ws.Range("D16") '0.0000153'
$ws.Range("E16").Value = '  +2.83%  '
Set-LiteralText $ws.Range("D17") '2.564.01'
$ws.Range("E17").Value = '  -3.83%  '
Set-LiteralText $ws.Range("D18") '12.13'
$ws.Range("E18").Value = '  -0.59%  '
Set-LiteralText $ws.Range("D19") '4.66'
$ws.Range("E19").Value = '  -0.09%  '
Set-LiteralText $ws.Range("D20") '344.07'
$ws.Range("E20").Value = '  -1.22%  '
Set-LiteralText $ws.Range("D21") '6.79'
$ws.Range("E21").Value = '  -2.44%  '
Set-LiteralText $ws.Range("D22") '1.00'
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("E23").Value = '  +0.09%  '
Set-LiteralText $ws.Range("D24") '1.68'
$ws.Range("E24").Value = '  -3.65%  '
Set-LiteralText $ws.Range("D25") '9.06'
$ws.Range("E25").Value = '  -3.68%  '
Set-LiteralText $ws.Range("D26") '1.64'
$ws.Range("E26").Value = '  -3.90%  '
Set-LiteralText $ws.Range("D27") '553.55'
$ws.Range("E27").Value = '  -1.06%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-LiteralText $ws.Range("D28") '0.161'
$ws.Range("E28").Value = '  -1.99%  '
$ws.Range("B29").Value = 'Aptos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-LiteralText $ws.Range("D29") '7.97'
$ws.Range("E29").Value = '  -3.21%  '
$ws.Range("E30").Value = '  +0.39%  '
Set-LiteralText $ws.Range("D31") '2.01'
$ws.Range("E31").Value = '  -2.50%  '
Set-LiteralText $ws.Range("D32") '0.0₃0851'
$ws.Range("E32").Value = '  -2.06%  '
$ws.Range("E33").Value = '  -2.23%  '
Set-LiteralText $ws.Range("D34") '5.16'
$ws.Range("E34").Value = '  -4.11%  '
Set-LiteralText $ws.Range("D35") '165.15'
$ws.Range("E35").Value = '  -2.06%  '
Set-LiteralText $ws.Range("D36") '0.412'
$ws.Range("E36").Value = '  +0.74%  '
Set-LiteralText $ws.Range("D37") '0.999'
$ws.Range("E37").Value = '  -0.19%  '
Set-LiteralText $ws.Range("D38") '19.33'
$ws.Range("E38").Value = '  -0.60%  '
Set-LiteralText $ws.Range("D39") '1.88'
$ws.Range("E39").Value = '  -4.88%  '
$ws.Range("E40").Value = '  +0.04%  '
Set-LiteralText $ws.Range("D41") '165.10'
$ws.Range("E41").Value = '  -1.42%  '
Set-LiteralText $ws.Range("D42") '39.64'
$ws.Range("E42").Value = '  -1.79%  '
Set-LiteralText $ws.Range("D43") '3.94'
$ws.Range("E43").Value = '  +1.94%  '
Set-LiteralText $ws.Range("D44") '0.0586'
$ws.Range("E44").Value = '  +1.35%  '
Set-LiteralText $ws.Range("D45") '22.49'
$ws.Range("E45").Value = '  +1.58%  '
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-LiteralText $ws.Range("D46") '0.626'
$ws.Range("E46").Value = '  -1.02%  '
$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-LiteralText $ws.Range("D47") '2.03'
$ws.Range("E47").Value = '  +1.08%  '
$ws.Range("E48").Value = '  +0.18%  '
Set-LiteralText $ws.Range("D49") '0.0959'
$ws.Range("E49").Value = '  -0.89%  '
Set-LiteralText $ws.Range("D50") '18.86'
$ws.Range("E50").Value = '  -1.48%  '
Set-LiteralText $ws.Range("D51") '0.0₆0226'
$ws.Range("E51").Value = '  +11.13%  '
